$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new quarterly columns at D:E. This shifts the existing D:K data
# (quarters) right to F:M.
# ---------------------------------------------------------------------------
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy the number/date formatting from the (now shifted) old "D" column,
# which landed in F, into the two freshly inserted D:E columns so that the
# new cells carry the same styling (date format for header row, number
# format for data rows) as the rest of the table.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("G7:G102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# New data for the two inserted quarters (column D = newest quarter,
# column E = quarter before it), keyed by worksheet row number.
# ---------------------------------------------------------------------------
$newData = @{
    7 = @(43465, 43373)
    8 = @(436900, 445500)
    9 = @(331800, 326700)
    10 = @(105100, 118800)
    11 = @($null, $null)
    12 = @(1200, 1700)
    13 = @(0, 0)
    14 = @(-1700, 1600)
    15 = @(0, 0)
    16 = @($null, $null)
    17 = @(393200, 388500)
    18 = @(43700, 57000)
    19 = @($null, $null)
    20 = @(1300, 1500)
    21 = @(77800, 91900)
    22 = @(8900, 9700)
    23 = @(36100, 48800)
    24 = @(4300, 13800)
    25 = @(0, 0)
    26 = @(31800, 35000)
    27 = @(30000, 33200)
    28 = @(0, 0)
    29 = @(15900, -300)
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(-1300, -1500)
    33 = @(45900, 32800)
    34 = @(0, 0)
    35 = @(45900, 32800)
    38 = @(43465, 43373)
    39 = @($null, $null)
    40 = @($null, $null)
    41 = @(64300, 61700)
    42 = @(0, 0)
    43 = @(369600, 387200)
    44 = @(133100, 137800)
    45 = @(38000, 44900)
    46 = @(605000, 631600)
    47 = @("NA", 3600)
    48 = @(469900, 460500)
    49 = @(491400, 503900)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(66600, 60800)
    53 = @(0, 0)
    54 = @(1632900, 1660400)
    55 = @($null, $null)
    56 = @($null, $null)
    57 = @(149400, 149200)
    58 = @(16600, 14800)
    59 = @(251000, 267100)
    60 = @(417000, 431100)
    61 = @(585700, 625400)
    62 = @(316800, 299500)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(1364600, 1398900)
    67 = @($null, $null)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(1298800, 1252800)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(268300, 261600)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(45900, 32800)
    82 = @($null, $null)
    83 = @(32800, 33400)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(97000, 48300)
    90 = @($null, $null)
    91 = @(-40900, -34800)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-25800, -22700)
    95 = @($null, $null)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-68900, -27300)
    101 = @(200, -900)
    102 = @(2600, -2500)
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    if ($vals[0] -ne $null) {
        $ws.Cells.Item($row, 4).Value = $vals[0]
    }
    if ($vals[1] -ne $null) {
        $ws.Cells.Item($row, 5).Value = $vals[1]
    }
}

# ---------------------------------------------------------------------------
# A handful of historical quarters (now in columns H and I) were corrected
# at the same time the new quarters were added.
# ---------------------------------------------------------------------------
$corrections = @{
    9  = @(350400, 290300)
    10 = @(104600, 94400)
    17 = @(418200, 350000)
    18 = @(36800, 34700)
    20 = @(300, -100)
    32 = @(-300, 100)
}

foreach ($row in $corrections.Keys) {
    $vals = $corrections[$row]
    $ws.Cells.Item($row, 8).Value = $vals[0]
    $ws.Cells.Item($row, 9).Value = $vals[1]
}
